$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "ECs" is introduced as a cluster label alongside the
# existing "FAPs" / "sCs" labels (commit: "Natmi following Dr Hou advice").
# The data table is rebuilt to cover all combinations of the three sending
# clusters (ECs, FAPs, sCs) against the single ligand/receptor pair
# (Clcf1/Crlf1) and three target clusters (ECs, FAPs, sCs), expanding the
# sheet from 4 data rows to 9 data rows (rows 2-10).

$rows = @(
    @(2, "ECs", "Clcf1", "Crlf1", "ECs", 1, 0.3333333333333333, 0.9214586666666666, 2.764376, 0.08041853843186561, 0.08041853843186561, 2, 0.6666666666666666, 0.1676186666666667, 0.502856, 0.005977118127032088, 0.005977118127032088, 0.1544536730951111, 1.390083057856, 0.0004806711038105305, 0.0004806711038105305),
    @(3, "ECs", "Clcf1", "Crlf1", "FAPs", 1, 0.3333333333333333, 0.9214586666666666, 2.764376, 0.08041853843186561, 0.08041853843186561, 3, 1, 16.17191166666667, 48.51573500000001, 0.5766745929546137, 0.5766745929546137, 14.90174816181778, 134.11573345636, 0.04637532791620105, 0.04637532791620105),
    @(4, "ECs", "Clcf1", "Crlf1", "sCs", 1, 0.3333333333333333, 0.9214586666666666, 2.764376, 0.08041853843186561, 0.08041853843186561, 3, 1, 11.70386166666667, 35.111585, 0.4173482889183543, 0.4173482889183543, 10.78462476621778, 97.06162289596, 0.03356253941185403, 0.03356253941185403),
    @(5, "FAPs", "Clcf1", "Crlf1", "ECs", 3, 1, 1.913147, 5.739441, 0.1669662363715809, 0.1669662363715809, 2, 0.6666666666666666, 0.1676186666666667, 0.502856, 0.005977118127032088, 0.005977118127032088, 0.3206791492773334, 2.886112343496, 0.0009979769180189003, 0.0009979769180189003),
    @(6, "FAPs", "Clcf1", "Crlf1", "FAPs", 3, 1, 1.913147, 5.739441, 0.1669662363715809, 0.1669662363715809, 3, 1, 16.17191166666667, 48.51573500000001, 0.5766745929546137, 0.5766745929546137, 30.93924428934834, 278.453198604135, 0.0962851863967452, 0.0962851863967452),
    @(7, "FAPs", "Clcf1", "Crlf1", "sCs", 3, 1, 1.913147, 5.739441, 0.1669662363715809, 0.1669662363715809, 3, 1, 11.70386166666667, 35.111585, 0.4173482889183543, 0.4173482889183543, 22.39120783599833, 201.520870523985, 0.06968307305681677, 0.06968307305681677),
    @(8, "sCs", "Clcf1", "Crlf1", "ECs", 3, 1, 8.623680999999999, 25.871043, 0.7526152251965536, 0.7526152251965536, 2, 0.6666666666666666, 0.1676186666666667, 0.502856, 0.005977118127032088, 0.005977118127032088, 1.445489910978667, 13.009409198808, 0.004498470105202657, 0.004498470105202657),
    @(9, "sCs", "Clcf1", "Crlf1", "FAPs", 3, 1, 8.623680999999999, 25.871043, 0.7526152251965536, 0.7526152251965536, 3, 1, 16.17191166666667, 48.51573500000001, 0.5766745929546137, 0.5766745929546137, 139.4614073735117, 1255.152666361605, 0.4340140786416674, 0.4340140786416674),
    @(10, "sCs", "Clcf1", "Crlf1", "sCs", 3, 1, 8.623680999999999, 25.871043, 0.7526152251965536, 0.7526152251965536, 3, 1, 11.70386166666667, 35.111585, 0.4173482889183543, 0.4173482889183543, 100.9303694814617, 908.373325333155, 0.3141026764496835, 0.3141026764496835)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rows) {
    $r = $row[0]
    for ($i = 1; $i -lt $row.Length; $i++) {
        $col = $cols[$i - 1]
        $ws.Range("$col$r").Value = $row[$i]
    }
}
